$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the dSF column (F) values for the rows whose underlying data was
# repulled/recalculated, per the commit: "repull data, push all data, mean calculation"
$ws.Range("F3").Value = 0
$ws.Range("F4").Value = -5
$ws.Range("F5").Value = -4
$ws.Range("F6").Value = -8
$ws.Range("F8").Value = -9
$ws.Range("F15").Value = -2
